$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# ---------------------------------------------------------------------------
# 1. Update the ConceptScheme URI (prefix base changes from '#' to '/')
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "http://data.sparna.fr/vocabularies/days/"

# ---------------------------------------------------------------------------
# 2. Insert a new explanatory row right after row 8 (pushes the header/table
#    down by one row) and fill it in with the new editorial note sentence.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert()
$ws.Range("A8").Font.Bold = $false
$ws.Range("A9").Font.Bold = $false
$ws.Range("A9").Value = "It also illustrates how the deprecattion of old concepts could be managed with owl:deprecated + dct:replacedBy"

# ---------------------------------------------------------------------------
# 3. New columns F/G/H on the header row (row 11 after the insert) plus the
#    widened/new column widths.
# ---------------------------------------------------------------------------
$ws.Range("D11").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F11").Value = "owl:deprecated^^xsd:boolean"
$ws.Range("G11").Value = "skos:editorialNote@en"
$ws.Range("H11").Value = "dct:replacedBy"

$ws.Columns.Item(6).ColumnWidth = 25.45
$ws.Columns.Item(7).ColumnWidth = 38.65
$ws.Columns.Item(8).ColumnWidth = 22.1

# ---------------------------------------------------------------------------
# 4. New row 19: the deprecated "days:primidi" example concept.
# ---------------------------------------------------------------------------
$ws.Range("C18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Copy() | Out-Null
$ws.Range("G19").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Copy() | Out-Null
$ws.Range("H19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").Value = "days:primidi"
$ws.Range("C19").Value = "Primidi"
$ws.Range("F19").Value = "true"
$ws.Range("G19").Value = "This was the first day of a `"decade`" (10 days) under the french revolution. See https://en.wikipedia.org/wiki/French_Republican_calendar#Ten_days_of_the_week`nLet's consider for the sake of the example that this is a deprecated concept that should be replaced by `"Monday`""
$ws.Range("H19").Value = "days:monday"

$ws.Rows.Item(19).RowHeight = 120

$ws.Range("A9").Select() | Out-Null
